# Refatorando o consolidador para modelo ETL
# Atualiza os registros de absenteísmo (linhas 2-11) com os novos dados
# vindos do pipeline ETL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Array de novos dados: ID, Nome, Departamento, Motivo, Horas, Data(serial), Salario
$data = @(
    @(19929, "Carlos Eduardo Costa", "Recursos Humanos", "Consulta médica", 2, 45097, 7939.36),
    @(97013, "Valentina Rocha", "P&D", "Consulta médica", 7, 45098, 9315.26),
    @(14771, "Laís Pereira", "Recursos Humanos", "Problemas pessoais", 3, 45081, 9443.75),
    @(13868, "Leonardo Fernandes", "Operações", "Outros", 3, 45094, 10783.82),
    @(97243, "Joaquim da Costa", "Financeiro", "Doença", 1, 45102, 8126.23),
    @(65333, "Marcela Ribeiro", "Engenharia", "Consulta médica", 2, 45093, 3142.56),
    @(50278, "Luiz Otávio Rodrigues", "Engenharia", "Viagem de negócios", 7, 45106, 5481.49),
    @(12654, "Emilly Ribeiro", "Jurídico", "Viagem de negócios", 3, 45102, 10347.09),
    @(94680, "Nina Mendes", "Vendas", "Consulta médica", 6, 45098, 10372.02),
    @(84723, "Leandro Caldeira", "Operações", "Problemas pessoais", 7, 45103, 11510.32)
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
    $ws.Cells.Item($row, 5).Value = $record[4]
    $ws.Cells.Item($row, 6).Value = $record[5]
    $ws.Cells.Item($row, 7).Value = $record[6]
    $row++
}
